# Apply the "想去人数" (want-to-go count) refresh + a couple of content
# corrections to 北京-漫展信息.xlsx, matching the upstream gh-pages
# regeneration at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")      # sheet1
$wsShow = $wb.Worksheets.Item("演出")      # sheet2
$wsAll  = $wb.Worksheets.Item("全部类型")  # sheet4

# ---------------------------------------------------------------------
# 展览 (Exhibitions) — column F ("想去人数") bumps
# ---------------------------------------------------------------------
$expoUpdates = @(
    @(2, 248), @(4, 9461), @(5, 611), @(7, 286), @(8, 353),
    @(9, 409), @(11, 185), @(13, 442), @(14, 12089), @(24, 2727),
    @(25, 2100), @(30, 1007), @(31, 4202), @(32, 3645), @(33, 556),
    @(41, 108), @(42, 424), @(43, 532), @(44, 69), @(48, 128)
)
foreach ($pair in $expoUpdates) {
    $row = $pair[0]
    $newVal = $pair[1]
    $wsExpo.Range("F$row").Value2 = $newVal
}

# ---------------------------------------------------------------------
# 演出 (Performances) — column F bumps
# ---------------------------------------------------------------------
$wsShow.Range("F3").Value2 = 6
$wsShow.Range("F9").Value2 = 43
$wsShow.Range("F17").Value2 = 15

# G6 (最低票价 for the Disney concert) became unavailable for sale
$wsShow.Range("G6").Value2 = "不可售"

# Rows 12/13/14 (all dated 2024-05-12) were re-ordered upstream; the
# A (index) and B (date) columns stay put, but the remaining columns
# C..I rotate: new12 <- old13, new13 <- old14, new14 <- old12.
$cols = @("C", "D", "E", "F", "G", "H", "I")

$orig12 = @{}
$orig13 = @{}
$orig14 = @{}
foreach ($c in $cols) {
    $orig12[$c] = $wsShow.Range("$c" + "12").Value2
    $orig13[$c] = $wsShow.Range("$c" + "13").Value2
    $orig14[$c] = $wsShow.Range("$c" + "14").Value2
}

foreach ($c in $cols) {
    $wsShow.Range("$c" + "12").Value2 = $orig13[$c]
    $wsShow.Range("$c" + "13").Value2 = $orig14[$c]
    $wsShow.Range("$c" + "14").Value2 = $orig12[$c]
}

# ---------------------------------------------------------------------
# 全部类型 (All types, merged view) — column F bumps
# ---------------------------------------------------------------------
$allUpdates = @(
    @(4, 248), @(6, 9461), @(7, 611), @(8, 43), @(10, 286),
    @(11, 353), @(12, 409), @(14, 185), @(15, 442), @(16, 12089),
    @(25, 2727), @(26, 2100), @(31, 1007), @(32, 4202), @(33, 3645),
    @(34, 556), @(40, 108), @(41, 424), @(43, 532), @(44, 69), @(48, 128)
)
foreach ($pair in $allUpdates) {
    $row = $pair[0]
    $newVal = $pair[1]
    $wsAll.Range("F$row").Value2 = $newVal
}
